$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header for column G (pop_sq_mile_1mi -> pop_sq_mile_10mi)
$ws.Range("G1").Value = "pop_sq_mile_10mi"

# Update total_risk (R) and total_risk_resp (S) values for rows 2-11
$ws.Range("R2").Value = 37.4242424242424
$ws.Range("S2").Value = 0.440909090909091

$ws.Range("R3").Value = 52.972972972973
$ws.Range("S3").Value = 0.510810810810811

$ws.Range("R4").Value = 29.9013157894737
$ws.Range("S4").Value = 0.3

$ws.Range("R5").Value = 27.1428571428571
$ws.Range("S5").Value = 0.314285714285714

$ws.Range("R6").Value = 30.2826086956522
$ws.Range("S6").Value = 0.370217391304348

$ws.Range("R7").Value = 29.2978208232446
$ws.Range("S7").Value = 0.362953995157385

$ws.Range("R8").Value = 53.9130434782609
$ws.Range("S8").Value = 0.427536231884058

$ws.Range("R9").Value = 53.8235294117647
$ws.Range("S9").Value = 0.463235294117647

$ws.Range("R10").Value = 20
$ws.Range("S10").Value = 0.2

$ws.Range("R11").Value = 18.7878787878788
$ws.Range("S11").Value = 0.181818181818182
